$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# F8 text content was swapped to the Kiabi logo URL (same text already used in F13)
$ws.Range("F8").Value = "https://leparcduluc.fr/wp-content/uploads/2022/08/1200px-Kiabi_logo.svg.png"

# New plain (non-hyperlinked) URL cells
$ws.Range("F3").Value = "https://fr.packcity.com/Packcity/images/logoPckCity.png"
$ws.Range("F11").Value = "https://upload.wikimedia.org/wikipedia/fr/thumb/d/d6/Logo_Worldline_-_2021.svg/800px-Logo_Worldline_-_2021.svg.png"

# F6: brand new URL text, with a self-referencing hyperlink (auto-applies the Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("F6"), "https://static.reserved.com/media/SHARED/stronywizerunkowe/reserved/cms/help/klarna_pink.png")

# F14: existing URL text now gets a self-referencing hyperlink (auto-applies the Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("F14"), "https://creation-entreprise.info/wp-content/uploads/2024/05/smart-rh.jpg")

# Update selection to match the final saved view state
$ws.Range("F11").Select()
